# Update automàtic: dades i banners [2026-02-06 02:49]
# Applies the refreshed MeteoCat daily summary values scraped at 2026-02-06 02:49 TU.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 02:47:49"
$ws.Range("N2").Value = "-2.0 °C 2:25 TU"
$ws.Range("O2").Value = "-0.9 °C"
$ws.Range("E3").Value = "2026-02-06 02:47:51"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "78%"
$ws.Range("O3").Value = "-2.1 °C"
$ws.Range("E4").Value = "2026-02-06 02:47:53"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "55%"
$ws.Range("J4").Value = "991.8 hPa"
$ws.Range("N4").Value = "12.7 °C 2:27 TU"
$ws.Range("O4").Value = "13.8 °C"
$ws.Range("E5").Value = "2026-02-06 02:47:55"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "71%"
$ws.Range("J5").Value = "992.5 hPa"
$ws.Range("L5").Value = "15.1 km/h - 298º 2:20 TU"
$ws.Range("N5").Value = "7.5 °C 2:08 TU"
$ws.Range("O5").Value = "9.0 °C"
$ws.Range("E6").Value = "2026-02-06 02:47:58"
$ws.Range("J6").Value = "993.6 hPa"
$ws.Range("N6").Value = "14.2 °C 2:29 TU"
$ws.Range("O6").Value = "14.8 °C"
$ws.Range("E7").Value = "2026-02-06 02:48:00"
$ws.Range("J7").Value = "993.4 hPa"
$ws.Range("N7").Value = "10.0 °C 2:09 TU"
$ws.Range("E8").Value = "2026-02-06 02:48:03"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "89%"
$ws.Range("N8").Value = "5.3 °C 2:18 TU"
$ws.Range("O8").Value = "6.8 °C"
$ws.Range("E9").Value = "2026-02-06 02:48:05"
$ws.Range("O9").Value = "2.7 °C"
$ws.Range("E10").Value = "2026-02-06 02:48:07"
$ws.Range("N10").Value = "4.8 °C 2:11 TU"
$ws.Range("O10").Value = "5.7 °C"
$ws.Range("E11").Value = "2026-02-06 02:48:10"
$ws.Range("J11").Value = "994.5 hPa"
$ws.Range("E12").Value = "2026-02-06 02:48:12"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "56%"
$ws.Range("O12").Value = "13.4 °C"
$ws.Range("E13").Value = "2026-02-06 02:48:15"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "85%"
$ws.Range("N13").Value = "5.2 °C 2:21 TU"
$ws.Range("O13").Value = "7.6 °C"
$ws.Range("E14").Value = "2026-02-06 02:48:17"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "73%"
$ws.Range("N14").Value = "-3.9 °C 2:11 TU"
$ws.Range("O14").Value = "-3.4 °C"
$ws.Range("E15").Value = "2026-02-06 02:48:19"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "74%"
$ws.Range("J15").Value = "992.4 hPa"
$ws.Range("N15").Value = "6.4 °C 2:00 TU"
$ws.Range("O15").Value = "9.1 °C"
$ws.Range("E16").Value = "2026-02-06 02:48:22"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "93%"
$ws.Range("O16").Value = "4.4 °C"
$ws.Range("E17").Value = "2026-02-06 02:48:24"
$ws.Range("J17").Value = "996.1 hPa"
$ws.Range("E18").Value = "2026-02-06 02:48:27"
$ws.Range("N18").Value = "-4.8 °C 2:29 TU"
$ws.Range("O18").Value = "-4.6 °C"
$ws.Range("E19").Value = "2026-02-06 02:48:29"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "96%"
$ws.Range("J19").Value = "996.4 hPa"
$ws.Range("N19").Value = "4.7 °C 2:29 TU"
$ws.Range("O19").Value = "7.1 °C"
$ws.Range("E20").Value = "2026-02-06 02:48:31"
$ws.Range("N20").Value = "-3.2 °C 2:27 TU"
$ws.Range("O20").Value = "-1.5 °C"
$ws.Range("E21").Value = "2026-02-06 02:48:34"
$ws.Range("H21").NumberFormat = "@"
$ws.Range("H21").Value = "80%"
$ws.Range("J21").Value = "993.1 hPa"
$ws.Range("N21").Value = "3.7 °C 2:08 TU"
$ws.Range("O21").Value = "6.2 °C"
$ws.Range("E22").Value = "2026-02-06 02:48:36"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "72%"
$ws.Range("N22").Value = "6.7 °C 2:20 TU"
$ws.Range("O22").Value = "10.4 °C"
$ws.Range("E23").Value = "2026-02-06 02:48:38"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "94%"
$ws.Range("J23").Value = "992.9 hPa"
$ws.Range("L23").Value = "14.8 km/h - 35º 2:10 TU"
$ws.Range("O23").Value = "6.9 °C"
$ws.Range("E24").Value = "2026-02-06 02:48:41"
$ws.Range("H24").NumberFormat = "@"
$ws.Range("H24").Value = "64%"
$ws.Range("J24").Value = "991.8 hPa"
$ws.Range("O24").Value = "12.2 °C"
$ws.Range("E25").Value = "2026-02-06 02:48:43"
$ws.Range("J25").Value = "994.9 hPa"
$ws.Range("E26").Value = "2026-02-06 02:48:46"
$ws.Range("L26").Value = "24.5 km/h - 28º 2:00 TU"
$ws.Range("N26").Value = "-0.6 °C 2:29 TU"
$ws.Range("E27").Value = "2026-02-06 02:48:48"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H27").Value = "94%"
$ws.Range("J27").Value = "992.4 hPa"
$ws.Range("N27").Value = "6.8 °C 2:29 TU"
$ws.Range("O27").Value = "8.3 °C"
$ws.Range("E28").Value = "2026-02-06 02:48:51"
$ws.Range("H28").NumberFormat = "@"
$ws.Range("H28").Value = "84%"
$ws.Range("J28").Value = "994.7 hPa"
$ws.Range("N28").Value = "1.9 °C 2:29 TU"
$ws.Range("O28").Value = "4.5 °C"
$ws.Range("E29").Value = "2026-02-06 02:48:53"
$ws.Range("H29").NumberFormat = "@"
$ws.Range("H29").Value = "53%"
$ws.Range("N29").Value = "12.8 °C 2:29 TU"
$ws.Range("O29").Value = "13.9 °C"
$ws.Range("E30").Value = "2026-02-06 02:48:56"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = "75%"
$ws.Range("N30").Value = "-4.5 °C 2:08 TU"
$ws.Range("O30").Value = "-3.2 °C"
$ws.Range("E31").Value = "2026-02-06 02:48:58"
$ws.Range("J31").Value = "996.1 hPa"
$ws.Range("O31").Value = "5.1 °C"
$ws.Range("E32").Value = "2026-02-06 02:49:00"
$ws.Range("J32").Value = "993.8 hPa"
$ws.Range("N32").Value = "15.3 °C 2:29 TU"
$ws.Range("E33").Value = "2026-02-06 02:49:02"
$ws.Range("N33").Value = "6.2 °C 2:27 TU"
$ws.Range("O33").Value = "7.6 °C"
$ws.Range("E34").Value = "2026-02-06 02:49:05"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "70%"
$ws.Range("N34").Value = "9.2 °C 2:28 TU"
$ws.Range("O34").Value = "9.7 °C"
$ws.Range("E35").Value = "2026-02-06 02:49:07"
$ws.Range("N35").Value = "-3.1 °C 2:29 TU"
$ws.Range("E36").Value = "2026-02-06 02:49:10"
$ws.Range("J36").Value = "995.5 hPa"
$ws.Range("O36").Value = "12.3 °C"
